# Registro actualización 2025-10-15 22:31:47
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7, column B (Cedula) currently holds the value as a text string.
# Convert it to a genuine numeric value, as the cedula is a plain number.
$ws.Cells.Item(7, 2).Value = 1000127336

# Append a new negotiation log entry as row 8, mirroring the structure
# of the existing rows, with the new timestamp 2025-10-15 17:31:47.
$ws.Cells.Item(8, 1).Value  = "2025-10-15 17:31:47"

# Cedula (B8) must stay textual, like the original row 7 value before it
# was normalised to a number. The leading apostrophe forces text entry;
# resetting the style back to Normal afterwards drops the quote-prefix
# formatting flag so the cell matches a plain text cell.
$ws.Cells.Item(8, 2).Value  = "'1000127336"
$ws.Cells.Item(8, 2).Style  = "Normal"

$ws.Cells.Item(8, 3).Value  = "Paula"
$ws.Cells.Item(8, 4).Value  = "TARJETA DE CRÉDITO"
$ws.Cells.Item(8, 5).Value  = "****4376"
$ws.Cells.Item(8, 6).Value  = "REDIFERIDO SIN PAGO"
$ws.Cells.Item(8, 7).Value  = "36 cuotas"
$ws.Cells.Item(8, 8).Value  = "34.19.100.134"
$ws.Cells.Item(8, 9).Value  = "The Dalles"
$ws.Cells.Item(8, 10).Value = "Oregon"
$ws.Cells.Item(8, 11).Value = "United States"
$ws.Cells.Item(8, 12).Value = "2025-10-15 17:31:47"
$ws.Cells.Item(8, 13).Value = "****4376"
$ws.Cells.Item(8, 14).Value = "34.19.100.134"

# O8/P8 are present but blank (empty text) in the source row layout.
# A plain "" assignment is treated as a cell-clear, so force an empty
# text cell via the quote-prefix trick and strip the resulting style.
$ws.Range("O8:P8").Value = "'"
$ws.Range("O8:P8").Style = "Normal"
